$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 679; existing rows 679:712 shift down to 680:713
$ws.Rows.Item(679).Insert()

# Populate the new row 679 with the new market entry
$ws.Cells.Item(679, 1).Value = 4
$ws.Cells.Item(679, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(679, 3).Value = "Los Lagos"
$ws.Cells.Item(679, 4).Value = 45041
$ws.Cells.Item(679, 5).Value = 10
$ws.Cells.Item(679, 6).Value = 100112006
$ws.Cells.Item(679, 7).Value = "Repollo"
$ws.Cells.Item(679, 8).Value = "Crespo record"
$ws.Cells.Item(679, 9).Value = "Primera"
$ws.Cells.Item(679, 10).Value = 1000
$ws.Cells.Item(679, 11).Value = 1800
$ws.Cells.Item(679, 12).Value = 2000
$ws.Cells.Item(679, 13).Value = 1900
$ws.Cells.Item(679, 14).Value = "$/unidad"
$ws.Cells.Item(679, 15).Value = "Región Metropolitana"
$ws.Cells.Item(679, 16).Value = 1900
$ws.Cells.Item(679, 17).Value = 1
$ws.Cells.Item(679, 18).Value = "Hortaliza"
